$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 125129.875
$ws.Range("I11").Value = 125129.875
$ws.Range("K11").Value = 125129.875
$ws.Range("M11").Value = -124989.875

# Row 101
$ws.Range("H101").Value = 429.72726
$ws.Range("I101").Value = 448.625
$ws.Range("J101").Value = 379.33334
$ws.Range("K101").Value = 1345.875
$ws.Range("L101").Value = 1138.00002
$ws.Range("M101").Value = 276.125
$ws.Range("N101").Value = -4382.000019999999

# Row 123
$ws.Range("H123").Value = 22035.4
$ws.Range("J123").Value = 22035.4
$ws.Range("L123").Value = 22035.4
$ws.Range("N123").Value = -31835.4

# Row 141
$ws.Range("H141").Value = 1969.32
$ws.Range("I141").Value = 1783.091
$ws.Range("J141").Value = 3335
$ws.Range("K141").Value = 5349.272999999999
$ws.Range("L141").Value = 10005
$ws.Range("M141").Value = -169.2729999999992
$ws.Range("N141").Value = -20365

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 111234160
$ws.Range("I110").Value = 125138056
$ws.Range("J110").Value = 2990
$ws.Range("K110").Value = 125138056
$ws.Range("L110").Value = 2990
$ws.Range("M110").Value = -125136011
$ws.Range("N110").Value = -7080

# Row 113
$ws.Range("H113").Value = 35032.668
$ws.Range("J113").Value = 35032.668
$ws.Range("L113").Value = 35032.668
$ws.Range("N113").Value = -43710.668

# Row 124
$ws.Range("H124").Value = 27186.715
$ws.Range("J124").Value = 27186.715
$ws.Range("L124").Value = 27186.715
$ws.Range("N124").Value = -37006.715

# Row 125
$ws.Range("H125").Value = 38850
$ws.Range("J125").Value = 38850
$ws.Range("L125").Value = 38850
$ws.Range("N125").Value = -48690

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1571.5625
$ws.Range("I20").Value = 1852.1
$ws.Range("J20").Value = 1104
$ws.Range("K20").Value = 1852.1
$ws.Range("L20").Value = 1104
$ws.Range("M20").Value = -1605.1
$ws.Range("N20").Value = -1598

# Row 86
$ws.Range("H86").Value = 85923
$ws.Range("I86").Value = 138417.38
$ws.Range("J86").Value = 1932
$ws.Range("K86").Value = 138417.38
$ws.Range("L86").Value = 1932
$ws.Range("M86").Value = -137294.38
$ws.Range("N86").Value = -4178

# Row 89
$ws.Range("H89").Value = 85923
$ws.Range("I89").Value = 138417.38
$ws.Range("J89").Value = 1932
$ws.Range("K89").Value = 692086.9
$ws.Range("L89").Value = 9660
$ws.Range("M89").Value = -686470.9
$ws.Range("N89").Value = -20892

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 852
$ws.Range("I17").Value = 852
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 852
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -678
$ws.Range("N17").ClearContents()

# Row 31
$ws.Range("H31").Value = 17010.213
$ws.Range("I31").Value = 24278.07
$ws.Range("J31").Value = 3422.4783
$ws.Range("K31").Value = 24278.07
$ws.Range("L31").Value = 3422.4783
$ws.Range("M31").Value = -23983.07
$ws.Range("N31").Value = -4012.4783

# Row 34
$ws.Range("H34").Value = 17010.213
$ws.Range("I34").Value = 24278.07
$ws.Range("J34").Value = 3422.4783
$ws.Range("K34").Value = 24278.07
$ws.Range("L34").Value = 3422.4783
$ws.Range("M34").Value = -24076.07
$ws.Range("N34").Value = -3826.4783

# Row 41
$ws.Range("H41").Value = 10370
$ws.Range("I41").Value = 3812.5
$ws.Range("J41").Value = 12387.692
$ws.Range("K41").Value = 3812.5
$ws.Range("L41").Value = 12387.692
$ws.Range("M41").Value = -3384.5
$ws.Range("N41").Value = -13243.692

# Row 86
$ws.Range("H86").Value = 2083.111
$ws.Range("I86").Value = 1850
$ws.Range("J86").Value = 2269.6
$ws.Range("K86").Value = 1850
$ws.Range("L86").Value = 2269.6
$ws.Range("M86").Value = -727
$ws.Range("N86").Value = -4515.6

# Row 89
$ws.Range("H89").Value = 2083.111
$ws.Range("I89").Value = 1850
$ws.Range("J89").Value = 2269.6
$ws.Range("K89").Value = 9250
$ws.Range("L89").Value = 11348
$ws.Range("M89").Value = -3634
$ws.Range("N89").Value = -22580

# Row 132
$ws.Range("H132").Value = 75004170
$ws.Range("I132").Value = 83337780
$ws.Range("J132").Value = 62503756
$ws.Range("K132").Value = 250013340
$ws.Range("L132").Value = 187511268
$ws.Range("M132").Value = -250010810
$ws.Range("N132").Value = -187516328

# Row 134
$ws.Range("H134").Value = 1073.0454
$ws.Range("I134").Value = 1039.2778
$ws.Range("J134").Value = 1225
$ws.Range("K134").Value = 3117.8334
$ws.Range("L134").Value = 3675
$ws.Range("M134").Value = -582.8334000000004
$ws.Range("N134").Value = -8745

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 737290.6
$ws.Range("J37").Value = 737290.6
$ws.Range("L37").Value = 2211871.8
$ws.Range("N37").Value = -2212095.8

# Row 58
$ws.Range("H58").Value = 1073.75
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1073.75
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3221.25
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3477.25

# Row 92
$ws.Range("H92").Value = 284
$ws.Range("J92").Value = 285.66666
$ws.Range("L92").Value = 856.9999799999999
$ws.Range("N92").Value = -3352.99998

# Row 131
$ws.Range("H131").Value = 843.13
$ws.Range("I131").Value = 606.44446
$ws.Range("J131").Value = 866.53845
$ws.Range("K131").Value = 1819.33338
$ws.Range("L131").Value = 2599.61535
$ws.Range("M131").Value = 3220.66662
$ws.Range("N131").Value = -12679.61535

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 172119.17
$ws.Range("I70").Value = 255001.5
$ws.Range("J70").Value = 6354.5
$ws.Range("K70").Value = 255001.5
$ws.Range("L70").Value = 6354.5
$ws.Range("M70").Value = -254731.5
$ws.Range("N70").Value = -6894.5

# Row 73
$ws.Range("H73").Value = 172119.17
$ws.Range("I73").Value = 255001.5
$ws.Range("J73").Value = 6354.5
$ws.Range("K73").Value = 255001.5
$ws.Range("L73").Value = 6354.5
$ws.Range("M73").Value = -254065.5
$ws.Range("N73").Value = -8226.5

# Row 97
$ws.Range("H97").Value = 76925680
$ws.Range("I97").Value = 111114056
$ws.Range("J97").Value = 1833
$ws.Range("K97").Value = 111114056
$ws.Range("L97").Value = 1833
$ws.Range("M97").Value = -111113560
$ws.Range("N97").Value = -2825

# Row 107
$ws.Range("H107").Value = 350.86667
$ws.Range("I107").Value = 343
$ws.Range("J107").Value = 359.85715
$ws.Range("K107").Value = 343
$ws.Range("L107").Value = 359.85715
$ws.Range("M107").Value = 1577
$ws.Range("N107").Value = -4199.85715

# Row 109
$ws.Range("H109").Value = 7017.5
$ws.Range("J109").Value = 7017.5
$ws.Range("L109").Value = 7017.5
$ws.Range("N109").Value = -9097.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 921428.2
$ws.Range("I46").Value = 1250
$ws.Range("J46").Value = 1447244.2
$ws.Range("K46").Value = 1250
$ws.Range("L46").Value = 1447244.2
$ws.Range("M46").Value = -1062
$ws.Range("N46").Value = -1447620.2

# Row 61
$ws.Range("H61").Value = 2018.1923
$ws.Range("I61").Value = 1628.238
$ws.Range("J61").Value = 3656
$ws.Range("K61").Value = 1628.238
$ws.Range("L61").Value = 3656
$ws.Range("M61").Value = -1426.238
$ws.Range("N61").Value = -4060

# Row 113
$ws.Range("H113").Value = 2018.1923
$ws.Range("I113").Value = 1628.238
$ws.Range("J113").Value = 3656
$ws.Range("K113").Value = 1628.238
$ws.Range("L113").Value = 3656
$ws.Range("M113").Value = 541.7619999999999
$ws.Range("N113").Value = -7996

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1416
$ws.Range("I136").Value = 395.81818
$ws.Range("J136").Value = 2484.762
$ws.Range("K136").Value = 1187.45454
$ws.Range("L136").Value = 7454.286
$ws.Range("M136").Value = 1362.54546
$ws.Range("N136").Value = -12554.286
